$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. Insert a row above the
# current row 130 (shifting existing rows 130-151 down to 131-152) and
# populate it with the new reading.
$ws.Rows.Item(130).Insert()

$ws.Range("A130").Value = 11
$ws.Range("B130").Value = "Vega Monumental Concepción"
$ws.Range("C130").Value = "Bíobío"
$ws.Range("D130").Value = 44798
$ws.Range("E130").Value = 8
$ws.Range("F130").Value = 100112043
$ws.Range("G130").Value = "Pepino ensalada"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 100
$ws.Range("K130").Value = 20000
$ws.Range("L130").Value = 22000
$ws.Range("M130").Value = 21000
$ws.Range("N130").Value = "$/caja 60 unidades"
$ws.Range("O130").Value = "Región de Arica y Parinacota"
$ws.Range("P130").Value = 350
$ws.Range("Q130").Value = 60
$ws.Range("R130").Value = "Hortaliza"

# Match the source style of the date column used by the rest of the sheet.
$ws.Range("D130").Style = $ws.Range("D131").Style
